# regen sval data to filter save games
# Update the numeric stat columns (B:G) for rows 2-23 on the active sheet
# with the newly regenerated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.001754667048134761,0.3375848360084654,0.1529057820181812,0.4998867070740569,0,0.9921319921488383),
    @(3.182878228561681,1.65323645889881,0.1529057820181812,0.4998867070740569,1,5.488907176552729),
    @(3.182878228561681,1.65323645889881,0.1529057820181812,0.4998867070740569,0,5.488907176552729),
    @(3.182878228561681,1.65323645889881,0.7127328510149897,0.4998867070740569,1,6.048734245549538),
    @(1.505614041169197,1.65323645889881,0.7127328510149897,0.4998867070740569,1,4.371470058157054),
    @(0.7287194209349384,1.65323645889881,0.7127328510149897,0.4998867070740569,1,3.594575437922795),
    @(3.182878228561681,1.65323645889881,0.7127328510149897,0.4998867070740569,1,6.048734245549538),
    @(1.505614041169197,0.3375848360084654,0.1529057820181812,0.4998867070740569,0,2.495991366269901),
    @(3.182878228561681,1.65323645889881,3.082599426703578,0.4998867070740569,1,8.418600821238126),
    @(3.182878228561681,1.65323645889881,0.1529057820181812,0.4998867070740569,1,5.488907176552729),
    @(3.182878228561681,1.65323645889881,157.8057217802531,6.48142807727062,1,169.1232645449842),
    @(3.182878228561681,1.65323645889881,0.7127328510149897,0.4998867070740569,1,6.048734245549538),
    @(0.1554434735375247,0.004309184025731883,3.082599426703578,0.4998867070740569,0,3.742238791340892),
    @(3.182878228561681,1.65323645889881,0.1529057820181812,0.4998867070740569,1,5.488907176552729),
    @(1.505614041169197,1.65323645889881,16.98373111632243,0.4998867070740569,1,20.64246832346449),
    @(1.505614041169197,1.65323645889881,0.7127328510149897,0.4998867070740569,1,4.371470058157054),
    @(0.7287194209349384,0.05231270169004087,16.98373111632243,0.4998867070740569,1,18.26464994602146),
    @(1.505614041169197,1.65323645889881,16.98373111632243,0.4998867070740569,1,20.64246832346449),
    @(0.7287194209349384,1.65323645889881,0.1529057820181812,0.4998867070740569,0,3.034748368925986),
    @(1.505614041169197,1.65323645889881,0.1529057820181812,0.4998867070740569,1,3.811642989160245),
    @(3.182878228561681,1.65323645889881,0.7127328510149897,0.4998867070740569,1,6.048734245549538),
    @(1.505614041169197,0.3375848360084654,3.082599426703578,0.4998867070740569,1,5.425685010955299)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $col = $j + 2   # column B = 2
        $ws.Cells.Item($r, $col).Value = $rowVals[$j]
    }
}
